$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1985111662531018
$ws.Range("C2").Value = 0.533498759305211
$ws.Range("J2").Value = 0.0173697270471464
$ws.Range("O2").Value = 0.002481389578163772
$ws.Range("P2").Value = 0.1712158808933003
$ws.Range("S2").Value = 0.07692307692307693
$ws.Range("B3").Value = 0.01785714285714286
$ws.Range("C3").Value = 0.03125
$ws.Range("J3").Value = 0.03125
$ws.Range("P3").Value = 0.7633928571428571
$ws.Range("S3").Value = 0.15625
$ws.Range("J4").Value = 0.04081632653061224
$ws.Range("P4").Value = 0.5714285714285714
$ws.Range("S4").Value = 0.3877551020408163
$ws.Range("B6").Value = 0.06106870229007633
$ws.Range("D6").Value = 0.02290076335877863
$ws.Range("F6").Value = 0.09923664122137404
$ws.Range("J6").Value = 0.2557251908396946
$ws.Range("O6").Value = 0.04198473282442748
$ws.Range("Q6").Value = 0.1717557251908397
$ws.Range("R6").Value = 0.01908396946564886
$ws.Range("S6").Value = 0.3282442748091603
$ws.Range("B7").Value = 0.1214953271028037
$ws.Range("D7").Value = 0.02336448598130841
$ws.Range("E7").Value = 0.004672897196261682
$ws.Range("F7").Value = 0.04205607476635514
$ws.Range("J7").Value = 0.1495327102803738
$ws.Range("O7").Value = 0.01401869158878505
$ws.Range("Q7").Value = 0.2149532710280374
$ws.Range("R7").Value = 0.04672897196261682
$ws.Range("S7").Value = 0.3831775700934579
$ws.Range("B8").Value = 0.1036906854130053
$ws.Range("D8").Value = 0.01757469244288225
$ws.Range("F8").Value = 0.070298769771529
$ws.Range("J8").Value = 0.09314586994727592
$ws.Range("O8").Value = 0.01757469244288225
$ws.Range("Q8").Value = 0.210896309314587
$ws.Range("R8").Value = 0.08260105448154657
$ws.Range("S8").Value = 0.4042179261862918
$ws.Range("B9").Value = 0.1096491228070175
$ws.Range("D9").Value = 0.03508771929824561
$ws.Range("F9").Value = 0.07456140350877193
$ws.Range("J9").Value = 0.1271929824561404
$ws.Range("O9").Value = 0.01754385964912281
$ws.Range("Q9").Value = 0.206140350877193
$ws.Range("R9").Value = 0.03508771929824561
$ws.Range("S9").Value = 0.3947368421052632
$ws.Range("B10").Value = 0.1308864265927978
$ws.Range("D10").Value = 0.01731301939058172
$ws.Range("E10").Value = 0.001385041551246537
$ws.Range("F10").Value = 0.05955678670360111
$ws.Range("J10").Value = 0.103185595567867
$ws.Range("O10").Value = 0.01385041551246537
$ws.Range("Q10").Value = 0.2617728531855956
$ws.Range("R10").Value = 0.05540166204986149
$ws.Range("S10").Value = 0.3566481994459834
$ws.Range("F11").Value = 0.003144654088050315
$ws.Range("G11").Value = 0.1509433962264151
$ws.Range("J11").Value = 0.07547169811320754
$ws.Range("K11").Value = 0.1949685534591195
$ws.Range("L11").Value = 0.5628930817610063
$ws.Range("S11").Value = 0.01257861635220126
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.1904761904761905
$ws.Range("K12").Value = 0.005291005291005291
$ws.Range("L12").Value = 0.05291005291005291
$ws.Range("S12").Value = 0.03703703703703703
$ws.Range("F13").Value = 0.02272727272727273
$ws.Range("G13").Value = 0.7045454545454546
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.02272727272727273
$ws.Range("F15").Value = 0.03930131004366812
$ws.Range("H15").Value = 0.1615720524017467
$ws.Range("I15").Value = 0.1048034934497817
$ws.Range("J15").Value = 0.3624454148471616
$ws.Range("K15").Value = 0.06986899563318777
$ws.Range("M15").Value = 0.008733624454148471
$ws.Range("O15").Value = 0.0611353711790393
$ws.Range("S15").Value = 0.1921397379912664
$ws.Range("F16").Value = 0.01145038167938931
$ws.Range("H16").Value = 0.1679389312977099
$ws.Range("I16").Value = 0.08015267175572519
$ws.Range("J16").Value = 0.4160305343511451
$ws.Range("K16").Value = 0.09541984732824428
$ws.Range("M16").Value = 0.01908396946564886
$ws.Range("O16").Value = 0.03435114503816794
$ws.Range("S16").Value = 0.1755725190839695
$ws.Range("F17").Value = 0.01729559748427673
$ws.Range("H17").Value = 0.1886792452830189
$ws.Range("I17").Value = 0.0959119496855346
$ws.Range("J17").Value = 0.4292452830188679
$ws.Range("K17").Value = 0.0770440251572327
$ws.Range("M17").Value = 0.01729559748427673
$ws.Range("O17").Value = 0.07075471698113207
$ws.Range("S17").Value = 0.1037735849056604
$ws.Range("F18").Value = 0.02013422818791946
$ws.Range("H18").Value = 0.1812080536912752
$ws.Range("I18").Value = 0.08724832214765101
$ws.Range("J18").Value = 0.3624161073825503
$ws.Range("K18").Value = 0.1208053691275168
$ws.Range("M18").Value = 0.04026845637583892
$ws.Range("O18").Value = 0.04697986577181208
$ws.Range("S18").Value = 0.1409395973154362
$ws.Range("F19").Value = 0.01885474860335196
$ws.Range("H19").Value = 0.2374301675977654
$ws.Range("I19").Value = 0.07751396648044692
$ws.Range("J19").Value = 0.3666201117318436
$ws.Range("K19").Value = 0.1026536312849162
$ws.Range("M19").Value = 0.01745810055865922
$ws.Range("O19").Value = 0.05237430167597765
$ws.Range("S19").Value = 0.1270949720670391
